# "feat: add 2022-Q4 data"
#
# The workbook tracks one "总计" (summary) sheet plus one sheet per quarter
# of fund-holding data. This adds a new "2022-Q4" quarter:
#   - the previous "2022-Q3" sheet's numbers are refreshed in place and the
#     sheet is renamed to "2022-Q4"
#   - a fresh "2022-Q3" sheet is created alongside it holding the original
#     (now-historical) 2022-Q3 figures, so that data is preserved
#   - "2021-Q1" is untouched
#   - the "总计" summary sheet gets a new row for 2022-Q4

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Duplicate the current "2022-Q3" sheet (with all its data/formatting
# intact) so the original figures survive under the "2022-Q3" name once the
# source sheet becomes "2022-Q4".
$q3Sheet.Copy($null, $q3Sheet)
$q3CopySheet = $wb.Worksheets.Item($q3Sheet.Index + 1)

$q3Sheet.Name = "2022-Q4"
$q3CopySheet.Name = "2022-Q3"

$q4Sheet = $q3Sheet

# --- Refresh the "2022-Q4" sheet with the new quarter's numbers. ---
# These columns hold numeric-looking figures stored as text (e.g. "16.03",
# "0.7855") in the source data, so force a Text format before writing the
# values to avoid Excel silently converting them to numbers (which would
# also drop meaningful trailing zeros like the one in "4.90").
$q4Sheet.Range("D2:G4").NumberFormat = "@"

$q4Sheet.Range("D2").Value = "16.03"
$q4Sheet.Range("E2").Value = "98.69"
$q4Sheet.Range("F2").Value = "4.90"
$q4Sheet.Range("G2").Value = "0.7855"

$q4Sheet.Range("D3").Value = "14.15"
$q4Sheet.Range("E3").Value = "75.21"
$q4Sheet.Range("F3").Value = "4.99"
$q4Sheet.Range("G3").Value = "0.7061"
$q4Sheet.Range("H3").Value = 1

$q4Sheet.Range("E4").Value = "75.21"
$q4Sheet.Range("F4").Value = "4.99"
$q4Sheet.Range("G4").Value = "0.0195"
$q4Sheet.Range("H4").Value = 1

# --- Update the "总计" summary sheet: push existing rows down one and add
#     a new top row for 2022-Q4. ---
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 1.4

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.51
